$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.687.86"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.636.60"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("E6").Value = "  +4.06%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.31"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.90%  "
$ws.Range("E11").Value = "  +3.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.865.40"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.643.80"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.10"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.65%  "
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.694.11"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.56%  "
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0746"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "220.62"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +9.43%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.32"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.46"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.56%  "
$ws.Range("E23").Value = "  +2.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.93"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.29"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.98%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +1.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.95"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.54"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.57%  "
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.34"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.22%  "
$ws.Range("E33").Value = "  +2.46%  "
$ws.Range("E34").Value = "  +1.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.40"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.217.34"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0173"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.812"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.29%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("E40").Value = "  +2.56%  "
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("E42").Value = "  +2.13%  "
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.774.76"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.44"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("E46").Value = "  +2.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.93"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.26%  "
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("E49").Value = "  +5.69%  "
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("E51").Value = "  +0.01%  "
